# Apply cryptos.xlsx data update (Tue Jan 23 05:42:43 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain decimal-looking price values to remain text (avoid numeric auto-conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update changed cell values
$ws.Range("D2").Value = "40.121.10"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "2.346.90"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "310.55"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "85.39"
$ws.Range("E6").Value = "  -4.86%  "
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "0.0813"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "30.13"
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "2.706.06"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "6.42"
$ws.Range("E14").Value = "  -4.58%  "
$ws.Range("D15").Value = "14.80"
$ws.Range("E15").Value = "  -6.81%  "
$ws.Range("D16").Value = "2.367.41"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "0.761"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "40.070.94"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "0.0₃0902"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").Value = "6.10"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "68.20"
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("D22").Value = "10.67"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D23").Value = "235.11"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  -5.21%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("D27").Value = "23.82"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "34.71"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "153.83"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "15.59"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "1.965.57"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("D45").Value = "17.55"
$ws.Range("E45").Value = "  -7.25%  "
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("E47").Value = "  -7.00%  "
$ws.Range("D48").Value = "2.567.71"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "93.25"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").Value = "70.53"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "50.20"
$ws.Range("E51").Value = "  -3.84%  "
